$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 11333.333
$ws.Range("J7").Value = 11333.333
$ws.Range("L7").Value = 11333.333
$ws.Range("N7").Value = -11557.333

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 11333.333
$ws.Range("J14").Value = 11333.333
$ws.Range("L14").Value = 11333.333
$ws.Range("N14").Value = -11715.333

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 776.2
$ws.Range("I38").Value = 70.25
$ws.Range("J38").Value = 3600
$ws.Range("K38").Value = 210.75
$ws.Range("L38").Value = 10800
$ws.Range("M38").Value = 161.25
$ws.Range("N38").Value = -11544

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1547.1428
$ws.Range("I58").Value = 310
$ws.Range("J58").Value = 3196.6667
$ws.Range("K58").Value = 930
$ws.Range("L58").Value = 9590.000100000001
$ws.Range("M58").Value = -780
$ws.Range("N58").Value = -9890.000100000001

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1199.8334
$ws.Range("I2").Value = 1016.3333
$ws.Range("K2").Value = 1016.3333
$ws.Range("M2").Value = -903.3333

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13317.676
$ws.Range("I32").Value = 12140.099
$ws.Range("J32").Value = 27252.334
$ws.Range("K32").Value = 12140.099
$ws.Range("L32").Value = 27252.334
$ws.Range("M32").Value = -11853.099
$ws.Range("N32").Value = -27826.334

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2193.5833
$ws.Range("I74").Value = 1970.5
$ws.Range("J74").Value = 2416.6667
$ws.Range("K74").Value = 1970.5
$ws.Range("L74").Value = 2416.6667
$ws.Range("M74").Value = -1096.5
$ws.Range("N74").Value = -4164.6667

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2193.5833
$ws.Range("I77").Value = 1970.5
$ws.Range("J77").Value = 2416.6667
$ws.Range("K77").Value = 9852.5
$ws.Range("L77").Value = 12083.3335
$ws.Range("M77").Value = -5484.5
$ws.Range("N77").Value = -20819.3335

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1199.8334
$ws.Range("I116").Value = 1016.3333
$ws.Range("K116").Value = 1016.3333
$ws.Range("M116").Value = 1277.6667

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1199.8334
$ws.Range("I3").Value = 1016.3333
$ws.Range("K3").Value = 1016.3333
$ws.Range("M3").Value = -902.3333

# BSM row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -8096

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1843.8572
$ws.Range("I10").Value = 2951.75
$ws.Range("J10").Value = 366.66666
$ws.Range("K10").Value = 2951.75
$ws.Range("L10").Value = 366.66666
$ws.Range("M10").Value = -2812.75
$ws.Range("N10").Value = -644.66666

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1185.9
$ws.Range("I16").Value = 976.5
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 976.5
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -689.5
$ws.Range("N16").Value = -2074

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2924.08
$ws.Range("I99").Value = 2763.2632
$ws.Range("J99").Value = 3433.3333
$ws.Range("K99").Value = 2763.2632
$ws.Range("L99").Value = 3433.3333
$ws.Range("M99").Value = -1265.2632
$ws.Range("N99").Value = -6429.3333

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1185.9
$ws.Range("I113").Value = 976.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 976.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1193.5
$ws.Range("N113").Value = -5840

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2924.08
$ws.Range("I126").Value = 2763.2632
$ws.Range("J126").Value = 3433.3333
$ws.Range("K126").Value = 8289.7896
$ws.Range("L126").Value = 10299.9999
$ws.Range("M126").Value = -5819.7896
$ws.Range("N126").Value = -15239.9999

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1429.8928
$ws.Range("I134").Value = 1383.238
$ws.Range("J134").Value = 1569.8572
$ws.Range("K134").Value = 4149.714
$ws.Range("L134").Value = 4709.571599999999
$ws.Range("M134").Value = -1614.714
$ws.Range("N134").Value = -9779.571599999999

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1207777.4
$ws.Range("I12").Value = 32.75
$ws.Range("J12").Value = 1610358.9
$ws.Range("K12").Value = 98.25
$ws.Range("L12").Value = 4831076.699999999
$ws.Range("M12").Value = 74.75
$ws.Range("N12").Value = -4831422.699999999

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 256.2857
$ws.Range("I17").Value = 288
$ws.Range("J17").Value = 232.5
$ws.Range("K17").Value = 864
$ws.Range("L17").Value = 697.5
$ws.Range("M17").Value = -695
$ws.Range("N17").Value = -1035.5

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 814.6445
$ws.Range("J122").Value = 1015
$ws.Range("L122").Value = 9135
$ws.Range("N122").Value = -14035

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12502704
$ws.Range("J131").Value = 14707532
$ws.Range("L131").Value = 44122596
$ws.Range("N131").Value = -44132676

# GSM row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# GSM row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# GSM row 128
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 20900
$ws.Range("J128").Value = 20900
$ws.Range("L128").Value = 20900
$ws.Range("N128").Value = -30860

# GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 43754.285
$ws.Range("J140").Value = 43754.285
$ws.Range("L140").Value = 43754.285
$ws.Range("N140").Value = -54114.285

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 307.23077
$ws.Range("I22").Value = 339.4
$ws.Range("K22").Value = 339.4
$ws.Range("M22").Value = -44.39999999999998

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 307.23077
$ws.Range("I27").Value = 339.4
$ws.Range("K27").Value = 339.4
$ws.Range("M27").Value = -232.4

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8610.6875
$ws.Range("I40").Value = 7814.25
$ws.Range("J40").Value = 11000
$ws.Range("K40").Value = 7814.25
$ws.Range("L40").Value = 11000
$ws.Range("M40").Value = -7678.25
$ws.Range("N40").Value = -11272

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 18497.166
$ws.Range("I61").Value = 19724.182
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 19724.182
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -19522.182
$ws.Range("N61").Value = -5404

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1711.1111
$ws.Range("I93").Value = 1725
$ws.Range("J93").Value = 1700
$ws.Range("K93").Value = 1725
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = -477
$ws.Range("N93").Value = -4196

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 18497.166
$ws.Range("I113").Value = 19724.182
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 19724.182
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -17554.182
$ws.Range("N113").Value = -9340

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6364.7856
$ws.Range("I132").Value = 6396.6924
$ws.Range("J132").Value = 5950
$ws.Range("K132").Value = 19190.0772
$ws.Range("L132").Value = 17850
$ws.Range("M132").Value = -16660.0772
$ws.Range("N132").Value = -22910

# WVR row 9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 26751.75
$ws.Range("I9").Value = 26500
$ws.Range("J9").Value = 27003.5
$ws.Range("K9").Value = 26500
$ws.Range("L9").Value = 27003.5
$ws.Range("M9").Value = -26360
$ws.Range("N9").Value = -27283.5

# WVR row 87
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 143993.75
$ws.Range("J87").Value = 25325
$ws.Range("L87").Value = 25325
$ws.Range("N87").Value = -27821

# WVR row 90
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 143993.75
$ws.Range("J90").Value = 25325
$ws.Range("L90").Value = 75975
$ws.Range("N90").Value = -88455

# WVR row 92
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 32775
$ws.Range("J92").Value = 32775
$ws.Range("L92").Value = 32775
$ws.Range("N92").Value = -37767

# WVR row 93
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 38550
$ws.Range("J133").Value = 38550
$ws.Range("L133").Value = 38550
$ws.Range("N133").Value = -48670
